$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Site 3 (column D) representative coordinates corrected
$ws.Range("D3").Value2 = "42.4, -90.85"
$ws.Range("D4").Value2 = 42.4
$ws.Range("D5").Value2 = -90.85

# Site 5 (column F) representative coordinates corrected
$ws.Range("F3").Value2 = "47.5, -92.55"
$ws.Range("F4").Value2 = 47.5
$ws.Range("F5").Value2 = -92.55

# Update the active selection to K22 on Sheet1
$ws.Activate()
$ws.Range("K22").Select()
